$p = $ppt.ActivePresentation
Write-Host "HasHandoutMaster:" $p.HasHandoutMaster
try {
    $hm = $p.HandoutMaster
    Write-Host "HandoutMaster:" $hm
    $t3 = $hm.Theme
    $tcs3 = $t3.ThemeColorScheme
    Write-Host "Handout theme dk2:" $tcs3.Item(3).RGB
} catch {
    Write-Host "EXC:" $_.Exception.Message
}
